$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.864.48"
$ws.Range("E2").Value = "  -2.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.411.85"
$ws.Range("E3").Value = "  -1.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.95"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.91"
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.482"
$ws.Range("E8").Value = "  +0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.04"
$ws.Range("E9").Value = "  +5.06%  "

$ws.Range("E10").Value = "  -2.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.413"
$ws.Range("E11").Value = "  +2.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.995.43"
$ws.Range("E12").Value = "  -1.72%  "

$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.24"
$ws.Range("E14").Value = "  -5.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.390.78"
$ws.Range("E15").Value = "  -2.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.882.66"
$ws.Range("E17").Value = "  -1.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.36"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  +0.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.98"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.92"
$ws.Range("E21").Value = "  -2.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.566"
$ws.Range("E22").Value = "  +1.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.82"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.561.92"
$ws.Range("E25").Value = "  -1.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000111"
$ws.Range("E26").Value = "  -4.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.180"
$ws.Range("E27").Value = "  +0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.59"
$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.90"
$ws.Range("E30").Value = "  -3.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.11"
$ws.Range("E31").Value = "  -0.77%  "

$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("E33").Value = "  -3.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.03"
$ws.Range("E34").Value = "  -2.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  +2.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.61"
$ws.Range("E36").Value = "  +2.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "169.36"
$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("E38").Value = "  -2.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "30.89"
$ws.Range("E39").Value = "  -3.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.446.36"
$ws.Range("E40").Value = "  -1.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0776"
$ws.Range("E41").Value = "  +2.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.48"
$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.775"
$ws.Range("E43").Value = "  -3.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.36"
$ws.Range("E44").Value = "  -2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.67"
$ws.Range("E45").Value = "  -3.41%  "

$ws.Range("E46").Value = "  -4.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.542.99"
$ws.Range("E47").Value = "  -2.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.90"
$ws.Range("E48").Value = "  +2.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.63"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").Value = "  -4.11%  "

$ws.Range("E51").Value = "  +0.20%  "
